$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44824
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 8000
$ws.Cells.Item(2, 12).Value = 9000
$ws.Cells.Item(2, 13).Value = 8500
$ws.Cells.Item(2, 16).Value = 340
$ws.Cells.Item(3, 4).Value = 44848
$ws.Cells.Item(3, 10).Value = 800
$ws.Cells.Item(3, 11).Value = 7000
$ws.Cells.Item(3, 12).Value = 8000
$ws.Cells.Item(3, 13).Value = 7500
$ws.Cells.Item(3, 16).Value = 300
$ws.Cells.Item(4, 4).Value = 44837
$ws.Cells.Item(4, 10).Value = 520
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 9000
$ws.Cells.Item(4, 13).Value = 8500
$ws.Cells.Item(4, 16).Value = 340
$ws.Cells.Item(5, 4).Value = 44386
$ws.Cells.Item(5, 10).Value = 500
$ws.Cells.Item(5, 11).Value = 11000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 13).Value = 11500
$ws.Cells.Item(5, 16).Value = 460
$ws.Cells.Item(6, 4).Value = 44855
$ws.Cells.Item(6, 10).Value = 540
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 13).Value = 7500
$ws.Cells.Item(6, 16).Value = 300
$ws.Cells.Item(7, 4).Value = 44881
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 7000
$ws.Cells.Item(7, 12).Value = 8000
$ws.Cells.Item(7, 13).Value = 7500
$ws.Cells.Item(7, 16).Value = 300
$ws.Cells.Item(8, 4).Value = 44825
$ws.Cells.Item(8, 10).Value = 440
$ws.Cells.Item(8, 11).Value = 8000
$ws.Cells.Item(8, 13).Value = 8500
$ws.Cells.Item(8, 16).Value = 340
$ws.Cells.Item(9, 4).Value = 44384
$ws.Cells.Item(9, 10).Value = 560
$ws.Cells.Item(9, 11).Value = 11500
$ws.Cells.Item(9, 12).Value = 12000
$ws.Cells.Item(9, 13).Value = 11750
$ws.Cells.Item(9, 16).Value = 470
$ws.Cells.Item(10, 4).Value = 44466
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 9500
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 13).Value = 9750
$ws.Cells.Item(10, 16).Value = 390
$ws.Cells.Item(11, 4).Value = 44714
$ws.Cells.Item(11, 11).Value = 14000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 14500
$ws.Cells.Item(11, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 16).Value = 580
$ws.Cells.Item(12, 4).Value = 44880
$ws.Cells.Item(12, 10).Value = 560
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 7500
$ws.Cells.Item(12, 16).Value = 300
$ws.Cells.Item(13, 4).Value = 44799
$ws.Cells.Item(13, 11).Value = 10000
$ws.Cells.Item(13, 12).Value = 11000
$ws.Cells.Item(13, 13).Value = 10500
$ws.Cells.Item(13, 16).Value = 420
$ws.Cells.Item(14, 4).Value = 44377
$ws.Cells.Item(14, 10).Value = 520
$ws.Cells.Item(14, 11).Value = 12500
$ws.Cells.Item(14, 12).Value = 13000
$ws.Cells.Item(14, 13).Value = 12750
$ws.Cells.Item(14, 16).Value = 510
$ws.Cells.Item(15, 4).Value = 44803
$ws.Cells.Item(15, 10).Value = 600
$ws.Cells.Item(15, 11).Value = 9500
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 9750
$ws.Cells.Item(15, 16).Value = 390
$ws.Cells.Item(16, 4).Value = 44846
$ws.Cells.Item(16, 10).Value = 600
$ws.Cells.Item(16, 11).Value = 7500
$ws.Cells.Item(16, 12).Value = 8000
$ws.Cells.Item(16, 13).Value = 7750
$ws.Cells.Item(16, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(16, 16).Value = 310
$ws.Cells.Item(17, 4).Value = 44484
$ws.Cells.Item(17, 10).Value = 400
$ws.Cells.Item(17, 11).Value = 9000
$ws.Cells.Item(17, 12).Value = 10000
$ws.Cells.Item(17, 13).Value = 9500
$ws.Cells.Item(17, 16).Value = 380
$ws.Cells.Item(18, 4).Value = 44372
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 14000
$ws.Cells.Item(18, 13).Value = 13500
$ws.Cells.Item(18, 16).Value = 540
$ws.Cells.Item(19, 4).Value = 44876
$ws.Cells.Item(19, 10).Value = 460
$ws.Cells.Item(19, 11).Value = 6000
$ws.Cells.Item(19, 12).Value = 7000
$ws.Cells.Item(19, 13).Value = 6500
$ws.Cells.Item(19, 16).Value = 260
$ws.Cells.Item(20, 4).Value = 44811
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 10000
$ws.Cells.Item(20, 12).Value = 10500
$ws.Cells.Item(20, 13).Value = 10250
$ws.Cells.Item(20, 16).Value = 410
$ws.Cells.Item(21, 4).Value = 44874
$ws.Cells.Item(21, 10).Value = 500
$ws.Cells.Item(21, 11).Value = 6000
$ws.Cells.Item(21, 12).Value = 7000
$ws.Cells.Item(21, 13).Value = 6500
$ws.Cells.Item(21, 16).Value = 260
$ws.Cells.Item(22, 4).Value = 44756
$ws.Cells.Item(22, 11).Value = 14000
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 13).Value = 14500
$ws.Cells.Item(22, 16).Value = 580
$ws.Cells.Item(23, 4).Value = 44473
$ws.Cells.Item(23, 10).Value = 500
$ws.Cells.Item(23, 11).Value = 8500
$ws.Cells.Item(23, 12).Value = 9000
$ws.Cells.Item(23, 13).Value = 8750
$ws.Cells.Item(23, 16).Value = 350
$ws.Cells.Item(24, 4).Value = 44446
$ws.Cells.Item(24, 11).Value = 11000
$ws.Cells.Item(24, 12).Value = 12000
$ws.Cells.Item(24, 13).Value = 11500
$ws.Cells.Item(24, 16).Value = 460
$ws.Cells.Item(25, 4).Value = 44797
$ws.Cells.Item(25, 10).Value = 1000
$ws.Cells.Item(25, 11).Value = 11000
$ws.Cells.Item(25, 12).Value = 12000
$ws.Cells.Item(25, 13).Value = 11500
$ws.Cells.Item(25, 16).Value = 460
$ws.Cells.Item(26, 4).Value = 44809
$ws.Cells.Item(26, 10).Value = 520
$ws.Cells.Item(27, 4).Value = 44816
$ws.Cells.Item(27, 10).Value = 600
$ws.Cells.Item(27, 11).Value = 9500
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = 9750
$ws.Cells.Item(27, 16).Value = 390
$ws.Cells.Item(28, 4).Value = 44798
$ws.Cells.Item(28, 10).Value = 400
$ws.Cells.Item(28, 11).Value = 10500
$ws.Cells.Item(28, 12).Value = 11000
$ws.Cells.Item(28, 13).Value = 10750
$ws.Cells.Item(28, 16).Value = 430
$ws.Cells.Item(29, 4).Value = 44316
$ws.Cells.Item(29, 10).Value = 300
$ws.Cells.Item(29, 11).Value = 16000
$ws.Cells.Item(29, 12).Value = 17000
$ws.Cells.Item(29, 13).Value = 16500
$ws.Cells.Item(29, 16).Value = 660
$ws.Cells.Item(30, 4).Value = 44883
$ws.Cells.Item(30, 10).Value = 380
$ws.Cells.Item(30, 11).Value = 7000
$ws.Cells.Item(30, 12).Value = 8000
$ws.Cells.Item(30, 13).Value = 7500
$ws.Cells.Item(30, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(30, 16).Value = 300
$ws.Cells.Item(31, 4).Value = 44817
$ws.Cells.Item(31, 10).Value = 440
$ws.Cells.Item(31, 11).Value = 9000
$ws.Cells.Item(31, 13).Value = 9500
$ws.Cells.Item(31, 16).Value = 380
$ws.Cells.Item(32, 4).Value = 44858
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 8000
$ws.Cells.Item(32, 13).Value = 7500
$ws.Cells.Item(32, 16).Value = 300
$ws.Cells.Item(33, 4).Value = 44873
$ws.Cells.Item(33, 10).Value = 540
$ws.Cells.Item(33, 11).Value = 6000
$ws.Cells.Item(33, 12).Value = 7000
$ws.Cells.Item(33, 13).Value = 6500
$ws.Cells.Item(33, 16).Value = 260
$ws.Cells.Item(34, 4).Value = 44721
$ws.Cells.Item(34, 11).Value = 14500
$ws.Cells.Item(34, 12).Value = 15000
$ws.Cells.Item(34, 13).Value = 14750
$ws.Cells.Item(34, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 16).Value = 590
$ws.Cells.Item(35, 4).Value = 44781
$ws.Cells.Item(35, 11).Value = 10000
$ws.Cells.Item(35, 12).Value = 11000
$ws.Cells.Item(35, 13).Value = 10500
$ws.Cells.Item(35, 16).Value = 420
$ws.Cells.Item(36, 4).Value = 44694
$ws.Cells.Item(36, 10).Value = 480
$ws.Cells.Item(36, 11).Value = 17500
$ws.Cells.Item(36, 12).Value = 18000
$ws.Cells.Item(36, 13).Value = 17750
$ws.Cells.Item(36, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(36, 16).Value = 710
$ws.Cells.Item(39, 4).Value = 44425
$ws.Cells.Item(39, 11).Value = 11500
$ws.Cells.Item(39, 12).Value = 12000
$ws.Cells.Item(39, 13).Value = 11750
$ws.Cells.Item(39, 16).Value = 470
$ws.Cells.Item(40, 4).Value = 44827
$ws.Cells.Item(40, 10).Value = 700
$ws.Cells.Item(41, 4).Value = 44690
$ws.Cells.Item(41, 10).Value = 400
$ws.Cells.Item(41, 11).Value = 17000
$ws.Cells.Item(41, 12).Value = 18000
$ws.Cells.Item(41, 13).Value = 17500
$ws.Cells.Item(41, 16).Value = 700
$ws.Cells.Item(42, 4).Value = 44370
$ws.Cells.Item(42, 10).Value = 520
$ws.Cells.Item(42, 11).Value = 13000
$ws.Cells.Item(42, 12).Value = 14000
$ws.Cells.Item(42, 13).Value = 13500
$ws.Cells.Item(42, 16).Value = 540
$ws.Cells.Item(43, 4).Value = 44356
$ws.Cells.Item(43, 10).Value = 500
$ws.Cells.Item(43, 11).Value = 13000
$ws.Cells.Item(43, 12).Value = 14000
$ws.Cells.Item(43, 13).Value = 13500
$ws.Cells.Item(43, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(43, 16).Value = 540
